$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1, styled like the other header cells (bold, centered)
$ws.Range("E1").Value = "Colocação"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108

# New column E values (ranking / "Colocação") for rows 2-7
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 6
